$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (8:40): add a "Hardwareentwicklung" class in Dienstag (C4) and
# an "Englisch" class in Donnerstag (E4)
$ws.Range("C4").Value = "Hardwareentwicklung`nGERE`nW11x"
$ws.Range("E4").Value = "Englisch`nRAKL`nW118"

# Row 6 (9:40): remove the "Mathematik" entry that used to be in Montag (B6)
$ws.Range("B6").Value = ""

# Row 7 (10:30): remove "Deutsch" from Montag (B7) and add "Mathematik" to
# Mittwoch (D7)
$ws.Range("B7").Value = ""
$ws.Range("D7").Value = "Mathematik`nWÜBE`nW118"

# Row 8 (11:20): remove the "Softwareentwicklung" entry from Montag (B8)
$ws.Range("B8").Value = ""

# Column widths changed because column B lost its long text and column D
# gained a longer one, so Excel recomputed the best-fit widths.
$ws.Columns.Item(2).ColumnWidth = 7.7265625
$ws.Columns.Item(4).ColumnWidth = 11.48828125
